# Replace the three occurrences of the old site URL ("goibibo.com") with the
# new one ("makemytrip.com"). Per the target diff, in every case the literal
# hostname text "makemytrip" ends up isolated in its own <w:r> run (sitting
# between a run with "https://www." and a run with ".com/..."), even though
# its resulting run-properties are identical to its neighbours' properties.
#
# Word's normal text-replace collapses same-formatted adjacent runs back
# into one, so to keep the post-edit run boundaries in place we briefly
# toggle a run-level formatting property on just the "makemytrip" span and
# then toggle it back to its original value. That round-trip is what keeps
# the span as a physically separate run after the edit.

$d = $word.ActiveDocument

function Replace-Host([string]$oldFullText, [string]$newFullText, [bool]$hostIsBold) {
    $full = $d.Content
    $null = $full.Find.Execute($oldFullText, $true, $false, $false, $false, $false, $true, 1, $false, $newFullText, 2)

    $hostRange = $d.Content
    $null = $hostRange.Find.Execute("makemytrip", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($hostIsBold) {
        # Already bold; dip to non-bold and back so the run survives as its
        # own element while ending on the original (bold) value.
        $hostRange.Font.Bold = 0
        $hostRange.Font.Bold = 1
    } else {
        # Already non-bold; flip to bold and back so the run survives as its
        # own element while ending on the original (non-bold) value.
        $hostRange.Font.Bold = 1
        $hostRange.Font.Bold = 0
    }
}

Replace-Host `
    "Test Strategy for Testing https://www.goibibo.com/" `
    "Test Strategy for Testing https://www.makemytrip.com/" `
    $true

Replace-Host `
    "   - The purpose of this test strategy is to outline the overall approach, objectives, and methodologies for testing the front-end functionality of the https://www.goibibo.com/ website." `
    "   - The purpose of this test strategy is to outline the overall approach, objectives, and methodologies for testing the front-end functionality of the https://www.makemytrip.com/ website." `
    $false

Replace-Host `
    "This test strategy provides a structured approach for testing the front-end functionality of the https://www.goibibo.com/ website, ensuring a high-quality user experience across different platforms and devices." `
    "This test strategy provides a structured approach for testing the front-end functionality of the https://www.makemytrip.com/ website, ensuring a high-quality user experience across different platforms and devices." `
    $false
